$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2159090909090909
$ws.Range("C2").Value = 0.5340909090909091
$ws.Range("J2").Value = 0.02556818181818182
$ws.Range("P2").Value = 0.15625
$ws.Range("S2").Value = 0.06818181818181818
$ws.Range("C3").Value = 0.01041666666666667
$ws.Range("J3").Value = 0.02083333333333333
$ws.Range("P3").Value = 0.71875
$ws.Range("S3").Value = 0.25
$ws.Range("J4").Value = 0.0851063829787234
$ws.Range("O4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.5957446808510638
$ws.Range("S4").Value = 0.2978723404255319
$ws.Range("B6").Value = 0.05314009661835749
$ws.Range("D6").Value = 0.004830917874396135
$ws.Range("F6").Value = 0.05797101449275362
$ws.Range("J6").Value = 0.3043478260869565
$ws.Range("O6").Value = 0.02415458937198068
$ws.Range("Q6").Value = 0.178743961352657
$ws.Range("R6").Value = 0.06763285024154589
$ws.Range("S6").Value = 0.3091787439613526
$ws.Range("B7").Value = 0.1023255813953488
$ws.Range("D7").Value = 0.0186046511627907
$ws.Range("F7").Value = 0.03255813953488372
$ws.Range("J7").Value = 0.1534883720930233
$ws.Range("O7").Value = 0.009302325581395349
$ws.Range("Q7").Value = 0.2046511627906977
$ws.Range("R7").Value = 0.1023255813953488
$ws.Range("S7").Value = 0.3767441860465116
$ws.Range("B8").Value = 0.1098901098901099
$ws.Range("D8").Value = 0.02417582417582418
$ws.Range("E8").Value = 0.002197802197802198
$ws.Range("F8").Value = 0.04615384615384616
$ws.Range("J8").Value = 0.1252747252747253
$ws.Range("O8").Value = 0.01978021978021978
$ws.Range("Q8").Value = 0.156043956043956
$ws.Range("R8").Value = 0.0989010989010989
$ws.Range("S8").Value = 0.4175824175824176
$ws.Range("B9").Value = 0.1194029850746269
$ws.Range("D9").Value = 0.01492537313432836
$ws.Range("F9").Value = 0.05472636815920398
$ws.Range("J9").Value = 0.06965174129353234
$ws.Range("O9").Value = 0.01492537313432836
$ws.Range("Q9").Value = 0.2537313432835821
$ws.Range("R9").Value = 0.07960199004975124
$ws.Range("S9").Value = 0.3930348258706468
$ws.Range("B10").Value = 0.1204111600587372
$ws.Range("D10").Value = 0.02349486049926578
$ws.Range("E10").Value = 0.0007342143906020558
$ws.Range("F10").Value = 0.05800293685756241
$ws.Range("J10").Value = 0.1380323054331865
$ws.Range("O10").Value = 0.01174743024963289
$ws.Range("Q10").Value = 0.2217327459618209
$ws.Range("R10").Value = 0.08370044052863436
$ws.Range("S10").Value = 0.342143906020558
$ws.Range("G11").Value = 0.1685082872928177
$ws.Range("J11").Value = 0.1104972375690608
$ws.Range("K11").Value = 0.2513812154696133
$ws.Range("L11").Value = 0.4419889502762431
$ws.Range("S11").Value = 0.02762430939226519
$ws.Range("G12").Value = 0.7751479289940828
$ws.Range("J12").Value = 0.1715976331360947
$ws.Range("K12").Value = 0.005917159763313609
$ws.Range("L12").Value = 0.02958579881656805
$ws.Range("S12").Value = 0.01775147928994083
$ws.Range("G13").Value = 0.6037735849056604
$ws.Range("J13").Value = 0.2641509433962264
$ws.Range("S13").Value = 0.1320754716981132
$ws.Range("F15").Value = 0.02521008403361345
$ws.Range("H15").Value = 0.1092436974789916
$ws.Range("I15").Value = 0.07983193277310924
$ws.Range("J15").Value = 0.3907563025210084
$ws.Range("K15").Value = 0.0546218487394958
$ws.Range("M15").Value = 0.01260504201680672
$ws.Range("O15").Value = 0.07142857142857142
$ws.Range("S15").Value = 0.2563025210084033
$ws.Range("F16").Value = 0.02777777777777778
$ws.Range("H16").Value = 0.125
$ws.Range("I16").Value = 0.09259259259259259
$ws.Range("J16").Value = 0.4074074074074074
$ws.Range("K16").Value = 0.1203703703703704
$ws.Range("M16").Value = 0.02777777777777778
$ws.Range("O16").Value = 0.07870370370370371
$ws.Range("S16").Value = 0.1203703703703704
$ws.Range("F17").Value = 0.01388888888888889
$ws.Range("H17").Value = 0.1646825396825397
$ws.Range("I17").Value = 0.125
$ws.Range("J17").Value = 0.4146825396825397
$ws.Range("K17").Value = 0.1031746031746032
$ws.Range("M17").Value = 0.01984126984126984
$ws.Range("N17").Value = 0.001984126984126984
$ws.Range("O17").Value = 0.06547619047619048
$ws.Range("S17").Value = 0.09126984126984126
$ws.Range("F18").Value = 0.04347826086956522
$ws.Range("H18").Value = 0.1256038647342995
$ws.Range("I18").Value = 0.07729468599033816
$ws.Range("J18").Value = 0.4202898550724637
$ws.Range("K18").Value = 0.0966183574879227
$ws.Range("M18").Value = 0.03381642512077294
$ws.Range("O18").Value = 0.05797101449275362
$ws.Range("S18").Value = 0.1449275362318841
$ws.Range("F19").Value = 0.02066772655007949
$ws.Range("H19").Value = 0.2329093799682035
$ws.Range("I19").Value = 0.06677265500794913
$ws.Range("J19").Value = 0.3569157392686804
$ws.Range("K19").Value = 0.1248012718600954
$ws.Range("M19").Value = 0.02305246422893482
$ws.Range("N19").Value = 0.000794912559618442
$ws.Range("O19").Value = 0.07233704292527822
$ws.Range("S19").Value = 0.1017488076311606
